$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'ECs'
$ws.Range('B2').Value = 'Ngf'
$ws.Range('C2').Value = 'Ngfr'
$ws.Range('D2').Value = 'FAPs'
$ws.Range('E2').Value = 1
$ws.Range('F2').Value = 0.3333333333333333
$ws.Range('G2').Value = 0.27169
$ws.Range('H2').Value = 0.81507
$ws.Range('I2').Value = 0.04814835840308388
$ws.Range('J2').Value = 0.04814835840308387
$ws.Range('K2').Value = 3
$ws.Range('L2').Value = 1
$ws.Range('M2').Value = 0.6946430000000001
$ws.Range('N2').Value = 2.083929
$ws.Range('O2').Value = 0.1140293552421611
$ws.Range('P2').Value = 0.1140293552421611
$ws.Range('Q2').Value = 0.18872755667
$ws.Range('R2').Value = 1.69854801003
$ws.Range('S2').Value = 0.005490326264672146
$ws.Range('T2').Value = 0.005490326264672146
# Row 3
$ws.Range('A3').Value = 'ECs'
$ws.Range('B3').Value = 'Ngf'
$ws.Range('C3').Value = 'Ngfr'
$ws.Range('D3').Value = 'MuSCs'
$ws.Range('E3').Value = 1
$ws.Range('F3').Value = 0.3333333333333333
$ws.Range('G3').Value = 0.27169
$ws.Range('H3').Value = 0.81507
$ws.Range('I3').Value = 0.04814835840308388
$ws.Range('J3').Value = 0.04814835840308387
$ws.Range('K3').Value = 3
$ws.Range('L3').Value = 1
$ws.Range('M3').Value = 5.373609333333333
$ws.Range('N3').Value = 16.120828
$ws.Range('O3').Value = 0.8821066470161785
$ws.Range('P3').Value = 0.8821066470161785
$ws.Range('Q3').Value = 1.459955919773333
$ws.Range('R3').Value = 13.13960327796
$ws.Range('S3').Value = 0.04247198699027756
$ws.Range('T3').Value = 0.04247198699027756
# Row 4
$ws.Range('A4').Value = 'ECs'
$ws.Range('B4').Value = 'Ngf'
$ws.Range('C4').Value = 'Ngfr'
$ws.Range('D4').Value = 'Resolving-Mac'
$ws.Range('E4').Value = 1
$ws.Range('F4').Value = 0.3333333333333333
$ws.Range('G4').Value = 0.27169
$ws.Range('H4').Value = 0.81507
$ws.Range('I4').Value = 0.04814835840308388
$ws.Range('J4').Value = 0.04814835840308387
$ws.Range('K4').Value = 1
$ws.Range('L4').Value = 0.3333333333333333
$ws.Range('M4').Value = 0.02353866666666667
$ws.Range('N4').Value = 0.070616
$ws.Range('O4').Value = 0.00386399774166032
$ws.Range('P4').Value = 0.00386399774166032
$ws.Range('Q4').Value = 0.006395220346666666
$ws.Range('R4').Value = 0.05755698312
$ws.Range('S4').Value = 0.0001860451481341678
$ws.Range('T4').Value = 0.0001860451481341678
# Row 5
$ws.Range('A5').Value = 'FAPs'
$ws.Range('B5').Value = 'Ngf'
$ws.Range('C5').Value = 'Ngfr'
$ws.Range('D5').Value = 'FAPs'
$ws.Range('E5').Value = 3
$ws.Range('F5').Value = 1
$ws.Range('G5').Value = 0.7246536666666668
$ws.Range('H5').Value = 2.173961
$ws.Range('I5').Value = 0.1284216734542145
$ws.Range('J5').Value = 0.1284216734542145
$ws.Range('K5').Value = 3
$ws.Range('L5').Value = 1
$ws.Range('M5').Value = 0.6946430000000001
$ws.Range('N5').Value = 2.083929
$ws.Range('O5').Value = 0.1140293552421611
$ws.Range('P5').Value = 0.1140293552421611
$ws.Range('Q5').Value = 0.5033755969743334
$ws.Range('R5').Value = 4.530380372769002
$ws.Range('S5').Value = 0.01464384062310344
$ws.Range('T5').Value = 0.01464384062310344
# Row 6
$ws.Range('A6').Value = 'FAPs'
$ws.Range('B6').Value = 'Ngf'
$ws.Range('C6').Value = 'Ngfr'
$ws.Range('D6').Value = 'MuSCs'
$ws.Range('E6').Value = 3
$ws.Range('F6').Value = 1
$ws.Range('G6').Value = 0.7246536666666668
$ws.Range('H6').Value = 2.173961
$ws.Range('I6').Value = 0.1284216734542145
$ws.Range('J6').Value = 0.1284216734542145
$ws.Range('K6').Value = 3
$ws.Range('L6').Value = 1
$ws.Range('M6').Value = 5.373609333333333
$ws.Range('N6').Value = 16.120828
$ws.Range('O6').Value = 0.8821066470161785
$ws.Range('P6').Value = 0.8821066470161785
$ws.Range('Q6').Value = 3.894005706634223
$ws.Range('R6').Value = 35.046051359708
$ws.Range('S6').Value = 0.1132816117749038
$ws.Range('T6').Value = 0.1132816117749038
# Row 7
$ws.Range('A7').Value = 'FAPs'
$ws.Range('B7').Value = 'Ngf'
$ws.Range('C7').Value = 'Ngfr'
$ws.Range('D7').Value = 'Resolving-Mac'
$ws.Range('E7').Value = 3
$ws.Range('F7').Value = 1
$ws.Range('G7').Value = 0.7246536666666668
$ws.Range('H7').Value = 2.173961
$ws.Range('I7').Value = 0.1284216734542145
$ws.Range('J7').Value = 0.1284216734542145
$ws.Range('K7').Value = 1
$ws.Range('L7').Value = 0.3333333333333333
$ws.Range('M7').Value = 0.02353866666666667
$ws.Range('N7').Value = 0.070616
$ws.Range('O7').Value = 0.00386399774166032
$ws.Range('P7').Value = 0.00386399774166032
$ws.Range('Q7').Value = 0.01705738110844444
$ws.Range('R7').Value = 0.153516429976
$ws.Range('S7').Value = 0.0004962210562073241
$ws.Range('T7').Value = 0.0004962210562073241
# Row 8
$ws.Range('A8').Value = 'MuSCs'
$ws.Range('B8').Value = 'Ngf'
$ws.Range('C8').Value = 'Ngfr'
$ws.Range('D8').Value = 'FAPs'
$ws.Range('E8').Value = 3
$ws.Range('F8').Value = 1
$ws.Range('G8').Value = 4.641773333333333
$ws.Range('H8').Value = 13.92532
$ws.Range('I8').Value = 0.8226057862976579
$ws.Range('J8').Value = 0.8226057862976578
$ws.Range('K8').Value = 3
$ws.Range('L8').Value = 1
$ws.Range('M8').Value = 0.6946430000000001
$ws.Range('N8').Value = 2.083929
$ws.Range('O8').Value = 0.1140293552421611
$ws.Range('P8').Value = 0.1140293552421611
$ws.Range('Q8').Value = 3.224375353586667
$ws.Range('R8').Value = 29.01937818228
$ws.Range('S8').Value = 0.09380120742999293
$ws.Range('T8').Value = 0.09380120742999291
# Row 9
$ws.Range('A9').Value = 'MuSCs'
$ws.Range('B9').Value = 'Ngf'
$ws.Range('C9').Value = 'Ngfr'
$ws.Range('D9').Value = 'MuSCs'
$ws.Range('E9').Value = 3
$ws.Range('F9').Value = 1
$ws.Range('G9').Value = 4.641773333333333
$ws.Range('H9').Value = 13.92532
$ws.Range('I9').Value = 0.8226057862976579
$ws.Range('J9').Value = 0.8226057862976578
$ws.Range('K9').Value = 3
$ws.Range('L9').Value = 1
$ws.Range('M9').Value = 5.373609333333333
$ws.Range('N9').Value = 16.120828
$ws.Range('O9').Value = 0.8821066470161785
$ws.Range('P9').Value = 0.8821066470161785
$ws.Range('Q9').Value = 24.94307650721778
$ws.Range('R9').Value = 224.48768856496
$ws.Range('S9').Value = 0.725626031967134
$ws.Range('T9').Value = 0.7256260319671339
# Row 10
$ws.Range('A10').Value = 'MuSCs'
$ws.Range('B10').Value = 'Ngf'
$ws.Range('C10').Value = 'Ngfr'
$ws.Range('D10').Value = 'Resolving-Mac'
$ws.Range('E10').Value = 3
$ws.Range('F10').Value = 1
$ws.Range('G10').Value = 4.641773333333333
$ws.Range('H10').Value = 13.92532
$ws.Range('I10').Value = 0.8226057862976579
$ws.Range('J10').Value = 0.8226057862976578
$ws.Range('K10').Value = 1
$ws.Range('L10').Value = 0.3333333333333333
$ws.Range('M10').Value = 0.02353866666666667
$ws.Range('N10').Value = 0.070616
$ws.Range('O10').Value = 0.00386399774166032
$ws.Range('P10').Value = 0.00386399774166032
$ws.Range('Q10').Value = 0.1092611552355555
$ws.Range('R10').Value = 0.9833503971199999
$ws.Range('S10').Value = 0.003178546900530862
$ws.Range('T10').Value = 0.003178546900530862
# Row 11
$ws.Range('A11').Value = 'Resolving-Mac'
$ws.Range('B11').Value = 'Ngf'
$ws.Range('C11').Value = 'Ngfr'
$ws.Range('D11').Value = 'FAPs'
$ws.Range('E11').Value = 1
$ws.Range('F11').Value = 0.3333333333333333
$ws.Range('G11').Value = 0.004650666666666667
$ws.Range('H11').Value = 0.013952
$ws.Range('I11').Value = 0.0008241818450437709
$ws.Range('J11').Value = 0.0008241818450437709
$ws.Range('K11').Value = 3
$ws.Range('L11').Value = 1
$ws.Range('M11').Value = 0.6946430000000001
$ws.Range('N11').Value = 2.083929
$ws.Range('O11').Value = 0.1140293552421611
$ws.Range('P11').Value = 0.1140293552421611
$ws.Range('Q11').Value = 0.003230553045333334
$ws.Range('R11').Value = 0.02907497740800001
$ws.Range('S11').Value = 0.00009398092439263596
$ws.Range('T11').Value = 0.00009398092439263596
# Row 12
$ws.Range('A12').Value = 'Resolving-Mac'
$ws.Range('B12').Value = 'Ngf'
$ws.Range('C12').Value = 'Ngfr'
$ws.Range('D12').Value = 'MuSCs'
$ws.Range('E12').Value = 1
$ws.Range('F12').Value = 0.3333333333333333
$ws.Range('G12').Value = 0.004650666666666667
$ws.Range('H12').Value = 0.013952
$ws.Range('I12').Value = 0.0008241818450437709
$ws.Range('J12').Value = 0.0008241818450437709
$ws.Range('K12').Value = 3
$ws.Range('L12').Value = 1
$ws.Range('M12').Value = 5.373609333333333
$ws.Range('N12').Value = 16.120828
$ws.Range('O12').Value = 0.8821066470161785
$ws.Range('P12').Value = 0.8821066470161785
$ws.Range('Q12').Value = 0.02499086580622222
$ws.Range('R12').Value = 0.224917792256
$ws.Range('S12').Value = 0.0007270162838631682
$ws.Range('T12').Value = 0.0007270162838631682
# Row 13
$ws.Range('A13').Value = 'Resolving-Mac'
$ws.Range('B13').Value = 'Ngf'
$ws.Range('C13').Value = 'Ngfr'
$ws.Range('D13').Value = 'Resolving-Mac'
$ws.Range('E13').Value = 1
$ws.Range('F13').Value = 0.3333333333333333
$ws.Range('G13').Value = 0.004650666666666667
$ws.Range('H13').Value = 0.013952
$ws.Range('I13').Value = 0.0008241818450437709
$ws.Range('J13').Value = 0.0008241818450437709
$ws.Range('K13').Value = 1
$ws.Range('L13').Value = 0.3333333333333333
$ws.Range('M13').Value = 0.02353866666666667
$ws.Range('N13').Value = 0.070616
$ws.Range('O13').Value = 0.00386399774166032
$ws.Range('P13').Value = 0.00386399774166032
$ws.Range('Q13').Value = 0.0001094704924444444
$ws.Range('R13').Value = 0.000985234432
$ws.Range('S13').Value = 0.000003184636787966566
$ws.Range('T13').Value = 0.000003184636787966566
